$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 42608.893379629626

$ws.Range("B6").Value = 85
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = "Random"
